# game-info.xlsx update: add new "ADDITIONAL" bonus column (I) for the
# ability-points table, bump up the DeadMines location numbers (one more
# elite kill + one more boss-item entry), and register the new
# "Tranquil Gardens" location under the additional-region notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "ADDITIONAL" bonus column header (I3) ------------------------
$ws.Range("I3").Value = "ADDITIONAL"

# --- Ability points table (rows 5-14): add flat "ADDITIONAL" bonus ----
# column I holds the flat bonus added on top of C*multiplier for rows 7-14.
$ws.Range("I7").Value = 5
$ws.Range("I8").Value = 7
$ws.Range("I9").Value = 9
$ws.Range("I10").Value = 11
$ws.Range("I11").Value = 13
$ws.Range("I12").Value = 15
$ws.Range("I13").Value = 17
$ws.Range("I14").Value = 19

$ws.Range("D7").Formula = '=C7*$D$4+I7'
$ws.Range("E7").Formula = '=C7*$E$4+I7'
$ws.Range("F7").Formula = '=C7*$F$4+I7'
$ws.Range("G7").Formula = '=C7*$G$4+I7'
$ws.Range("H7").Formula = '=C7*$H$4+I7'

$ws.Range("D8").Formula = '=C8*$D$4+I8'
$ws.Range("D9").Formula = '=C9*$D$4+I9'
$ws.Range("D10").Formula = '=C10*$D$4+I10'
$ws.Range("D11").Formula = '=C11*$D$4+I11'
$ws.Range("D12").Formula = '=C12*$D$4+I12'
$ws.Range("D13").Formula = '=C13*$D$4+I13'
$ws.Range("D14").Formula = '=C14*$D$4+I14'

$ws.Range("H8").Formula = '=C8*$H$4+I8'
$ws.Range("H9").Formula = '=C9*$H$4+I9'
$ws.Range("H10").Formula = '=C10*$H$4+I10'
$ws.Range("H11").Formula = '=C11*$H$4+I11'
$ws.Range("H12").Formula = '=C12*$H$4+I12'
$ws.Range("H13").Formula = '=C13*$H$4+I13'
$ws.Range("H14").Formula = '=C14*$H$4+I14'

# --- New region/location: DeadMines gains one more elite and loot slot,
#     bumping its elite-kill count (E21) and creature count (G21) -------
$ws.Range("E21").Value = 4
$ws.Range("G21").Value = 19

# --- New "Tranquil Gardens" location note ------------------------------
$ws.Range("E62").Value = "Tranquil Gardens"

# --- Restore default view (scroll position reset, new selection) ------
[void]$ws.Range("D15").Select()
